# FINFLUX-3612 Cartias specific scenarios
# Update Summary, Repayment schedule and Transactions sheets to reflect the
# revised (reduced) penalty-charge scenario figures, and move the active
# selection on several sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A5").Value = 0
$wsSummary.Range("B5").Value = 0
[void]$wsSummary.Range("C8").Select()

# ---------------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("J4").Value = 0
$wsRepay.Range("K4").Value = 42.94
$wsRepay.Range("L4").Value = 591.92999999999995
$wsRepay.Range("O4").Value = 42.94
[void]$wsRepay.Range("K7").Select()

# ---------------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------------
$wsTxn = $wb.Worksheets.Item("Transactions")

$wsTxn.Range("A2").Value = 664
$wsTxn.Range("A3").Value = 663
$wsTxn.Range("A4").Value = 662
$wsTxn.Range("E4").Value = 4533.79
$wsTxn.Range("F4").Value = 4493.95
$wsTxn.Range("A5").Value = 661
$wsTxn.Range("J5").Value = 4533.79
$wsTxn.Range("A6").Value = 660
$wsTxn.Range("A7").Value = 659
$wsTxn.Range("E7").Value = 43.67
$wsTxn.Range("I7").Value = 0
$wsTxn.Range("J7").Value = 4488.05
$wsTxn.Range("A8").Value = 658
$wsTxn.Range("E8").Value = 43.67
$wsTxn.Range("I8").Value = 0
$wsTxn.Range("A9").Value = 657
$wsTxn.Range("F9").Value = 506.05
$wsTxn.Range("I9").Value = 0
$wsTxn.Range("A10").Value = 652
$wsTxn.Range("A11").Value = 651
$wsTxn.Range("A12").Value = 650

[void]$wsTxn.Select()
[void]$wsTxn.Range("E5").Select()
